# Update the "weekly_source_path" value from the spring-promotion Weekly_DB
# folder to the new fuel-promotion Transactions folder, and relocate the
# (hidden) "_GoBack" bookmark so it sits right after the edited text instead
# of at the end of the values_to_skip paragraph.

$d = $word.ActiveDocument

$old = "Y:\_Current projects\Visa\Visa spring promotion 2019\Transactions\Weekly_DB"
$new = "Y:\_Current projects\Visa\Visa fuel promotion summer 2019\Transactions"

# Locate the run that holds the old path text.
$r = $d.Content
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not find the source text to replace."
}

$start = $r.Start

# Nudge the font size up for the duration of the edit so the rewritten run
# isn't silently coalesced into its (identically formatted) neighbouring
# runs; we restore the original size immediately afterwards.
$originalSize = $r.Font.Size
$r.Font.Size = $originalSize + 1
$r.Text = $new
$newEnd = $start + $new.Length

$r2 = $d.Range($start, $newEnd)
$r2.Font.Size = $originalSize

# Move the "_GoBack" bookmark to the end of the text we just inserted.
# Bookmarks.Add re-targets an existing bookmark of the same name rather than
# creating a duplicate, which also removes it from its old location (right
# after "TR_ID" in the values_to_skip paragraph).
$bmRange = $d.Range($newEnd, $newEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Replaced weekly_source_path value and moved _GoBack bookmark."
